$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 <- new values (previously row 6's data, with new date)
$ws.Range("D4").Value = 44334
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12500
$ws.Range("Q4").Value = '$/caja 12 kilos empedrada'
$ws.Range("S4").Value = 1042
$ws.Range("T4").Value = 12

# Row 5 <- new values (previously row 4's data)
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 17500
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17750
$ws.Range("S5").Value = 1109

# Row 6 <- new values (previously row 5's data, with new date)
$ws.Range("D6").Value = 44316
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = '$/caja 16 kilos granel'
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 16
